$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new rows of Mac-Address / machine data, following the existing pattern
$newRows = @(
    @{ A = 10001; B = 110030; C = 10030 },
    @{ A = 10001; B = 110031; C = 10031 }
)

$startRow = 31
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = "eng"
    $ws.Cells.Item($r, 5).Value = $true
    $ws.Cells.Item($r, 6).Value = "superadmin"
    $ws.Cells.Item($r, 7).Value = "now()"
    $ws.Cells.Item($r, 8).Value = "now()"
}

# Update the visible selection/scroll position like the final saved state
$ws.Application.ActiveWindow.ScrollRow = 25
$ws.Range("F30").Select()
